# Renames a handful of sector labels ("Power sector" -> "Power",
# "Land use" -> "LULUCF", and their derived/indented/pipe-joined variants)
# across the "Sheet 1" data sheet and the "labels" lookup sheet.

$wb = $excel.ActiveWorkbook

$data = $wb.Worksheets.Item("Sheet 1")
$labels = $wb.Worksheets.Item("labels")

# --- "Sheet 1": tier-2 sector name + tier-3 "Energy|..." path for the
#     Power sector row (code 1.A.1.a) and the LULUCF row (code 3B) ---
$data.Range("D2").Value = "Power"
$data.Range("G2").Value = "Energy|Power"

$data.Range("D33").Value = "LULUCF"
$data.Range("G33").Value = "AFOLU|LULUCF"

# --- "labels": indented display labels used for the legend/colour key ---
$labels.Range("A11").Value = "   LULUCF"
$labels.Range("A3").Value = "   Power"

# --- restore the last-selected cell on each sheet to match the saved file ---
$data.Range("G34").Select() | Out-Null
$labels.Range("A4").Select() | Out-Null
